$wb = $excel.ActiveWorkbook

# --- Final audit pass on "optimization_parameters": remove the stray
#     leftover "Sheet" row (row 16, values 3/4) that doesn't belong with
#     the rest of the optimization parameters table. Select the whole row
#     first (as a user would via the row header) then delete it; this
#     shifts the "simulation_timepoints" row up from 17 to 16.
$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Activate()
$wsOpt.Rows.Item(16).Select()
$wsOpt.Rows.Item(16).Delete()

# --- Move on to "network_weights" and point the selection at J27:K27
#     (anchored on K27) while reviewing the sheet.
$wsNet = $wb.Worksheets.Item("network_weights")
$wsNet.Activate()
$wsNet.Range("K27").Activate()
$wsNet.Range("J27:K27").Select()

# --- Finish the audit on "threshold_b", which ends up the active sheet
#     when the workbook is saved.
$wsThresh = $wb.Worksheets.Item("threshold_b")
$wsThresh.Activate()
$wsThresh.Range("A2").Select()
